$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the leftover "do nothing" test values in row 11 (A11 and C11),
# which also removes the now-unused shared strings "Be future" and "test".
$ws.Range("A11").ClearContents()
$ws.Range("C11").ClearContents()

# Update the selection to reflect the author's final cursor position.
$ws.Range("B18").Select()
